# Automatische test-sync: 2025-07-27 19:32:50
# Appends a new test-mail log entry (row 12) to the "Logs" sheet,
# extends the conditional-formatting ranges to cover it, and bumps the
# "Overig" tally on the "Dashboard" sheet from 4 to 5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 12

$ws.Cells.Item($row, 1).Value = "Is er al nieuws?"
$ws.Cells.Item($row, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value = "Testmail #10: Is er al nieuws?"
$ws.Cells.Item($row, 4).Value = "Overig"
$ws.Cells.Item($row, 5).Value = "Beste afzender,`nBedankt voor uw e-mail. Helaas kan ik u op basis van de informatie uit uw e-mail niet vertellen waar u precies naar op zoek bent. Als u meer context kunt geven of details over het onderwerp kunt verstrekken, help ik u graag verder. `nMet vriendelijke groet,`n[Je naam]`n[Bedrijfsnaam] E-mailassistent"
$ws.Cells.Item($row, 6).Value = "2025-07-27 19:32:24"
$ws.Cells.Item($row, 7).Value = "Ja"
$ws.Cells.Item($row, 8).Value = "Nee"
$ws.Cells.Item($row, 9).Value = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"

# Re-fit the row height so no stray explicit/custom row height sticks
# around after writing the multi-line "Antwoord" text (matches the
# other, untouched rows which have no ht/customHeight attribute).
$ws.Rows.Item($row).AutoFit()

# Extend the conditional formatting sqref for every coloured column
# from row 11 down to the newly added row 12.
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $col + "2:" + $col + "11"
    $newRange = $col + "2:" + $col + "12"
    $fc = $ws.Range($oldRange).FormatConditions
    $count = $fc.Count
    for ($i = 1; $i -le $count; $i++) {
        $fc.Item($i).ModifyAppliesToRange($ws.Range($newRange))
    }
}

# Dashboard: "Overig" count goes from 4 to 5.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 5
